$wb = $excel.ActiveWorkbook

# Sheet "展览": 熊喵M动漫嘉年华 (F3) and 万圣漫控嘉年华10 (F4)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 258
$wsExhibit.Range("F4").Value = 899

# Sheet "全部类型": 熊喵M动漫嘉年华 (F4) and 万圣漫控嘉年华10 (F5)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 258
$wsAll.Range("F5").Value = 899
